# Homework 10 edit — "Updated fee to 50cents and modified product backlog"
#
# 1. Product Backlog sheet: burndown "Left" slips by 2 more days starting
#    Sprint 2 (D10 formula gains "+2"), which cascades through B10:B13 and
#    D11:D13.
# 2. Product Backlog sheet: Features #7 and #8 (rows 27/28) are marked
#    Completed ("COM") with an Actual Sprint # of 2.
# 3. Sprint 2 sheet: the last two tasks finish a couple of days early
#    (B11 formula changes), and the two remaining "Not Started" backlog
#    rows (45, 48) are marked "COM".
# 4. Sprint 3 sheet gains three backlog rows (18-20) that were pulled in
#    from the Product Backlog.
# 5. Selection / active-sheet state moves from Sprint 2 back to the
#    Product Backlog sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1 & 2. Product Backlog sheet
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Product Backlog")

# Sprint 2 slips two extra days -> cascades to B10:B13 / D11:D13 via the
# existing formulas already on the sheet.
$backlog.Range("D10").Formula = "=D9+4+2"

# Features #7 (row 27) and #8 (row 28) are now Completed.
$backlog.Range("C27").Value = "COM"
$backlog.Range("E27").Value = 2
$backlog.Range("C28").Value = "COM"
$backlog.Range("E28").Value = 2

# ---------------------------------------------------------------------
# 3. Sprint 2 sheet
# ---------------------------------------------------------------------
$sprint2 = $wb.Worksheets.Item("Sprint 2")

$sprint2.Range("B11").Formula = "=B10 -3"

$sprint2.Range("G45").Value = "COM"
$sprint2.Range("G48").Value = "COM"

# ---------------------------------------------------------------------
# 4. Sprint 3 sheet gains three backlog rows
# ---------------------------------------------------------------------
$sprint3 = $wb.Worksheets.Item("Sprint 3")

$sprint3.Range("A18").Value = 1
$sprint3.Range("B18").Value = "Open a main GUI Window"

$sprint3.Range("A19").Formula = "=A18+1"
$sprint3.Range("B19").Value = "Create Media and Bundles Via GUI"

$sprint3.Range("B20").Value = "Browse GUI Catalog"

# ---------------------------------------------------------------------
# 5. Selection / active sheet moves back to Product Backlog
# ---------------------------------------------------------------------
$sprint2.Range("G45").Select()
$backlog.Select()
$backlog.Range("D11").Select()
